# A new weekly price record was added to the "Pomelo" sheet. It belongs
# right after the existing row 51, so row 52 onward (old rows 52:85) are
# pushed down by one to 53:86, and the freshly opened row 52 is filled in
# with the new record (same market/product classification as before,
# new date/volume/price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 52:85 down to 53:86 by inserting a blank row at 52.
$ws.Rows.Item(52).Insert()

# Fill in the new row 52 with the new record's data.
$ws.Cells.Item(52, 1).Value = 9                                              # A Mercado ID
$ws.Cells.Item(52, 2).Value = "Vega Central Mapocho de Santiago"             # B Mercado
$ws.Cells.Item(52, 3).Value = "Metropolitana"                                # C Región
$ws.Cells.Item(52, 4).Value = 44813                                          # D Fecha
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(53, 4).NumberFormat
$ws.Cells.Item(52, 5).Value = 13                                             # E Codreg
$ws.Cells.Item(52, 6).Value = "Fruta"                                        # F Tipo
$ws.Cells.Item(52, 7).Value = 100102                                         # G Producto ID
$ws.Cells.Item(52, 8).Value = "Cítricos"                                     # H Producto
$ws.Cells.Item(52, 9).Value = 100102006                                      # I Categoría ID
$ws.Cells.Item(52, 10).Value = "Pomelo"                                      # J Categoría
$ws.Cells.Item(52, 11).Value = "Start Ruby"                                  # K Variedad
$ws.Cells.Item(52, 12).Value = "Primera"                                     # L Calidad
$ws.Cells.Item(52, 13).Value = 300                                           # M Volumen
$ws.Cells.Item(52, 14).Value = 12000                                         # N Precio mínimo
$ws.Cells.Item(52, 15).Value = 12000                                         # O Precio máximo
$ws.Cells.Item(52, 16).Value = 12000                                         # P Precio promedio ponderado
$ws.Cells.Item(52, 17).Value = "$/caja 14 kilos"                             # Q Unidad de comercialización
$ws.Cells.Item(52, 18).Value = "Región Metropolitana"                        # R Origen
$ws.Cells.Item(52, 19).Value = 857                                           # S Precio $/Kg
$ws.Cells.Item(52, 20).Value = 14                                            # T Kg / unidad
